{"js": "const replacements = [\n  [\"60\u00d766=\", \"43\u00d781=\"],\n  [\"29\u00d741=\", \"61\u00d777=\"],\n  [\"44\u00d738=\", \"14\u00d761=\"],\n  [\"56\u00d749=\", \"26\u00d782=\"],\n  [\"82\u00d758=\", \"79\u00d715=\"],\n  [\"91\u00d757=\", \"51\u00d781=\"],\n  [\"96\u00d766=\", \"33\u00d730=\"],\n  [\"66\u00d758=\", \"64\u00d715=\"],\n  [\"13\u00d733=\", \"25\u00d751=\"],\n  [\"97\u00d798=\", \"81\u00d714=\"],\n  [\"37\u00d728=\", \"57\u00d770=\"],\n  [\"46\u00d765=\", \"20\u00d763=\"],\n  [\"43\u00d721=\", \"52\u00d786=\"],\n  [\"36\u00d714=\", \"19\u00d715=\"],\n  [\"54\u00d755=\", \"83\u00d763=\"],\n  [\"65\u00d774=\", \"91\u00d798=\"],\n  [\"74\u00d713=\", \"29\u00d742=\"],\n  [\"37\u00d776=\", \"39\u00d786=\"],\n  [\"69\u00d724=\", \"30\u00d751=\"],\n  [\"31\u00d749=\", \"41\u00d733=\"],\n  [\"63\u00d770=\", \"29\u00d775=\"],\n  [\"86\u00d773=\", \"74\u00d749=\"],\n  [\"46\u00d737=\", \"19\u00d783=\"],\n  [\"21\u00d769=\", \"68\u00d742=\"],\n  [\"85\u00d724=\", \"22\u00d773=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"60\u00d766=\"; New = \"43\u00d781=\" },\n    @{ Old = \"29\u00d741=\"; New = \"61\u00d777=\" },\n    @{ Old = \"44\u00d738=\"; New = \"14\u00d761=\" },\n    @{ Old = \"56\u00d749=\"; New = \"26\u00d782=\" },\n    @{ Old = \"82\u00d758=\"; New = \"79\u00d715=\" },\n    @{ Old = \"91\u00d757=\"; New = \"51\u00d781=\" },\n    @{ Old = \"96\u00d766=\"; New = \"33\u00d730=\" },\n    @{ Old = \"66\u00d758=\"; New = \"64\u00d715=\" },\n    @{ Old = \"13\u00d733=\"; New = \"25\u00d751=\" },\n    @{ Old = \"97\u00d798=\"; New = \"81\u00d714=\" },\n    @{ Old = \"37\u00d728=\"; New = \"57\u00d770=\" },\n    @{ Old = \"46\u00d765=\"; New = \"20\u00d763=\" },\n    @{ Old = \"43\u00d721=\"; New = \"52\u00d786=\" },\n    @{ Old = \"36\u00d714=\"; New = \"19\u00d715=\" },\n    @{ Old = \"54\u00d755=\"; New = \"83\u00d763=\" },\n    @{ Old = \"65\u00d774=\"; New = \"91\u00d798=\" },\n    @{ Old = \"74\u00d713=\"; New = \"29\u00d742=\" },\n    @{ Old = \"37\u00d776=\"; New = \"39\u00d786=\" },\n    @{ Old = \"69\u00d724=\"; New = \"30\u00d751=\" },\n    @{ Old = \"31\u00d749=\"; New = \"41\u00d733=\" },\n    @{ Old = \"63\u00d770=\"; New = \"29\u00d775=\" },\n    @{ Old = \"86\u00d773=\"; New = \"74\u00d749=\" },\n    @{ Old = \"46\u00d737=\"; New = \"19\u00d783=\" },\n    @{ Old = \"21\u00d769=\"; New = \"68\u00d742=\" },\n    @{ Old = \"85\u00d724=\"; New = \"22\u00d773=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
